$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update a few subcategory labels in column H
$ws.Range("H3").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H4").Value = "bar chart(s)"
$ws.Range("H7").Value = "drawing(s)"
$ws.Range("H8").Value = "line graph(s)"
$ws.Range("H14").Value = "line graph(s)"
$ws.Range("H19").Value = "line graph(s)"
$ws.Range("H20").Value = "line graph(s)"

# Remove the now-unwanted "is_viewed" column (column I) entirely
$ws.Columns.Item(9).Delete()
